$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 110, pushing existing rows 110:209 down to 111:210.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record.
$ws.Cells.Item(110, 1).Value = 11
$ws.Cells.Item(110, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(110, 3).Value = "Bíobío"
$ws.Cells.Item(110, 4).Value = 44586
$ws.Cells.Item(110, 5).Value = 8
$ws.Cells.Item(110, 6).Value = 100112009
$ws.Cells.Item(110, 7).Value = "Acelga"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 310
$ws.Cells.Item(110, 11).Value = 500
$ws.Cells.Item(110, 12).Value = 550
$ws.Cells.Item(110, 13).Value = 526
$ws.Cells.Item(110, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(110, 15).Value = "Región de Ñuble"
$ws.Cells.Item(110, 16).Value = 526
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"
